$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.841.01"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "3.743.54"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "404.46"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").Value = "127.36"
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("D7").Value = "3.732.18"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -6.02%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.716"
$ws.Range("E10").Value = "  -6.37%  "
$ws.Range("E11").Value = "  -9.32%  "
$ws.Range("D12").Value = "0.0000357"
$ws.Range("E12").Value = "  -6.33%  "
$ws.Range("D13").Value = "40.34"
$ws.Range("E13").Value = "  -5.56%  "
$ws.Range("D14").Value = "4.324.75"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "9.62"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").Value = "14.40"
$ws.Range("E16").Value = "  +11.60%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "3.721.32"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "19.40"
$ws.Range("E19").Value = "  -6.64%  "
$ws.Range("D20").Value = "66.144.12"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  -6.58%  "
$ws.Range("D22").Value = "405.48"
$ws.Range("E22").Value = "  -9.05%  "
$ws.Range("D23").Value = "14.41"
$ws.Range("E23").Value = "  -7.41%  "
$ws.Range("D24").Value = "84.86"
$ws.Range("E24").Value = "  -5.00%  "
$ws.Range("D25").Value = "3.00"
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("D26").Value = "36.08"
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("D27").Value = "5.55"
$ws.Range("E27").Value = "  +11.01%  "
$ws.Range("D28").Value = "3.08"
$ws.Range("E28").Value = "  -6.83%  "
$ws.Range("D29").Value = "9.25"
$ws.Range("E29").Value = "  -9.74%  "
$ws.Range("D30").Value = "12.35"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("D33").Value = "7.05"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("D35").Value = "38.38"
$ws.Range("E35").Value = "  -8.54%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "55.06"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").Value = "0.0₃0726"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "0.0453"
$ws.Range("E39").Value = "  -7.39%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  -8.87%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  -8.14%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "3.18"
$ws.Range("E43").Value = "  +21.84%  "
$ws.Range("D44").Value = "145.06"
$ws.Range("D45").Value = "26.46"
$ws.Range("E45").Value = "  -7.10%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "2.04"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "3.21"
$ws.Range("E47").Value = "  -7.31%  "
$ws.Range("D48").Value = "4.21"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").Value = "2.79"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("E51").Value = "  -6.54%  "
